$p = $ppt.ActivePresentation

function Merge-Runs($shape, [string]$finalText) {
    # Several runs in this paragraph carry the same visible sentence split
    # across multiple <a:r> elements (e.g. one run holds "embeddings" with
    # err="1" spell-check flag). The author's edit simply retyped the
    # sentence as one run. Re-assigning TextRange.Text directly is a no-op
    # here because the final string equals the existing concatenation, so
    # the host leaves the run boundaries untouched; deleting the old text
    # first (shrinking the shape) and then growing it back is also lossy
    # because autofit shapes don't reliably restore their original size.
    # Instead: append the desired final text at the end (only grows the
    # shape, never shrinks it), then delete the original prefix. What
    # remains is a single run using the first run's formatting, and the
    # shape's autofit size never dips below its original extent.
    $tr = $shape.TextFrame.TextRange
    $origLen = $tr.Length
    $tr.InsertAfter($finalText)
    $oldPrefix = $tr.Characters(1, $origLen)
    $oldPrefix.Delete()
}

# --- 1) Slide 12: fix typo "neithter" -> "Neithter" ---------------------
$s12 = $p.Slides.Item(12)
$sh12 = $s12.Shapes.Item(6)
$tr12 = $sh12.TextFrame.TextRange
$full12 = $tr12.Text
$idx12 = $full12.IndexOf("neithter")
if ($idx12 -ge 0) {
    $sub12 = $tr12.Characters($idx12 + 1, 8)
    $sub12.Text = "Neithter"
}

# --- 2) Slide 15: merge the "Visualization ... using PCA" runs ----------
$s15 = $p.Slides.Item(15)
$sh15 = $s15.Shapes.Item(7)
Merge-Runs $sh15 "Visualization of the Clusters for the document embeddings & project objective using PCA"

# --- 3) Slide 16: merge the "Elbow method ... claim embeddings." runs ---
$s16 = $p.Slides.Item(16)
$sh16a = $s16.Shapes.Item(6)
Merge-Runs $sh16a "Elbow method to choose the optimal number of clusters for the document embeddings & claim embeddings. "

# --- 4) Slide 16: merge the "Visualization ... claim embeddings using PCA" runs ---
$sh16b = $s16.Shapes.Item(7)
Merge-Runs $sh16b "Visualization of the Clusters for the document embeddings & claim embeddings using PCA"
